# Weekly update: a new price record (week of 2023-03-31) was added to the
# "Fruta, Terminal Hortofrutícola Agro Chillán - Pera" data set. The new
# record is inserted as row 172, which pushes every existing row from 172
# down to 263 (dimension grows from A1:T262 to A1:T263).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 172, shifting rows 172:262 down to 173:263.
$ws.Rows(172).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(172, 1).Value  = 7
$ws.Cells.Item(172, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(172, 3).Value  = "Ñuble"
$ws.Cells.Item(172, 4).Value  = 45016
$ws.Cells.Item(172, 5).Value  = 16
$ws.Cells.Item(172, 6).Value  = "Fruta"
$ws.Cells.Item(172, 7).Value  = 100104
$ws.Cells.Item(172, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(172, 9).Value  = 100104005
$ws.Cells.Item(172, 10).Value = "Pera"
$ws.Cells.Item(172, 11).Value = "Packham's Triumph"
$ws.Cells.Item(172, 12).Value = "Primera"
$ws.Cells.Item(172, 13).Value = 40
$ws.Cells.Item(172, 14).Value = 10000
$ws.Cells.Item(172, 15).Value = 10000
$ws.Cells.Item(172, 16).Value = 10000
$ws.Cells.Item(172, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(172, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(172, 19).Value = 556
$ws.Cells.Item(172, 20).Value = 18
